$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Standard Error (column D) values for rows 2-28, replacing the
# previous per-site constant SE with per-category SE values.
$newSE = @{
    2  = 0
    3  = 0.421878898054187
    4  = 0.749208680230307
    5  = 2.61953932570144
    6  = 0.178992957255875
    7  = 0.927469756457013
    8  = 0.138517774133618
    9  = 2.14258090758856
    10 = 0.627673996423881
    11 = 0.0031047367029617
    12 = 0.0402705308113329
    13 = 0.875309669693497
    14 = 1.43750216759478
    15 = 0.181282487758245
    16 = 0.73722057247343
    17 = 0.0705768819267477
    18 = 1.26128775176338
    19 = 0.286678173470934
    20 = 0.00289435600578871
    21 = 0.0565964355487878
    22 = 0.564762158499956
    23 = 1.09910804552703
    24 = 0.131172901144431
    25 = 0.835978576737014
    26 = 0.052411662439088
    27 = 1.00803275519988
    28 = 0.419306365653442
}

foreach ($row in $newSE.Keys) {
    $ws.Cells.Item($row, 4).Value = $newSE[$row]
}
